# Katalog guncellendi - Cum 21.11.2025 13:01:48,65
# Adds 6 new "FİYESMEN 1047 MONT" product rows (58-63) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 58: FİYESMEN 1047 MONT AÇIK YEŞİL ---------------------------------
$ws.Cells.Item(58, 1).Value = "FİYESMEN 1047 MONT AÇIK YEŞİL"
$ws.Cells.Item(58, 2).Value = "1300 TL"
$ws.Cells.Item(58, 3).Value = "Mont"
$ws.Cells.Item(58, 5).Value = "S-M-L-XL-2XL Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."
$ws.Cells.Item(58, 4).Value = "fiyesmen1047montaçıkyeşil.jpg"
$ws.Cells.Item(58, 6).Value = "Var"

# --- Row 59: FİYESMEN 1047 MONT ANTRASİT ------------------------------------
$ws.Cells.Item(59, 2).Value = "1300TL"
$ws.Cells.Item(59, 4).Value = "fiyesmen1047montantrasit.jpg"
$ws.Cells.Item(59, 1).Value = "FİYESMEN 1047 MONT ANTRASİT"
$ws.Cells.Item(59, 3).Value = "Mont"
$ws.Cells.Item(59, 5).Value = "S-M-L-XL-2XL Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."
$ws.Cells.Item(59, 6).Value = "Var"

# --- Row 60: FİYESMEN 1047 MONT GRİ -----------------------------------------
$ws.Cells.Item(60, 2).Value = "1300TL"
$ws.Cells.Item(60, 4).Value = "fiyesmen1047montgri.jpg"
$ws.Cells.Item(60, 1).Value = "FİYESMEN 1047 MONT GRİ"
$ws.Cells.Item(60, 3).Value = "Mont"
$ws.Cells.Item(60, 5).Value = "S-M-L-XL-2XL Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."
$ws.Cells.Item(60, 6).Value = "Var"

# --- Row 61: FİYESMEN 1047 MONT KAHVERENGİ ----------------------------------
$ws.Cells.Item(61, 2).Value = "1300TL"
$ws.Cells.Item(61, 4).Value = "fiyesmen1047montkahverengi.jpg"
$ws.Cells.Item(61, 1).Value = "FİYESMEN 1047 MONT KAHVERENGİ"
$ws.Cells.Item(61, 3).Value = "Mont"
$ws.Cells.Item(61, 5).Value = "S-M-L-XL-2XL Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."
$ws.Cells.Item(61, 6).Value = "Var"

# --- Row 62: FİYESMEN 1047 MONT SİYAH ---------------------------------------
$ws.Cells.Item(62, 2).Value = "1300TL"
$ws.Cells.Item(62, 4).Value = "fiyesmen1047montsiyah.jpg"
$ws.Cells.Item(62, 1).Value = "FİYESMEN 1047 MONT SİYAH"
$ws.Cells.Item(62, 3).Value = "Mont"
$ws.Cells.Item(62, 5).Value = "S-M-L-XL-2XL Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."
$ws.Cells.Item(62, 6).Value = "Var"

# --- Row 63: FİYESMEN 1047 MONT YEŞİL ---------------------------------------
$ws.Cells.Item(63, 2).Value = "1300TL"
$ws.Cells.Item(63, 4).Value = "fiyesmen1047montyeşil.jpg"
$ws.Cells.Item(63, 1).Value = "FİYESMEN 1047 MONT YEŞİL"
$ws.Cells.Item(63, 3).Value = "Mont"
$ws.Cells.Item(63, 5).Value = "S-M-L-XL-2XL Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."
$ws.Cells.Item(63, 6).Value = "Var"

# --- Column width adjustments (gorsel / aciklama columns widened) ----------
$ws.Columns.Item(4).ColumnWidth = 42.5
$ws.Columns.Item(5).ColumnWidth = 21.83

# --- View: scroll down and select F67 (matches author's final viewport) ----
$ws.Application.ActiveWindow.ScrollRow = 40
$ws.Range("F67").Select()
